$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Correct the D590:D603 timestamps (tiny precision fix from a re-save).
for ($r = 590; $r -le 603; $r++) {
    $ws.Cells.Item($r, 4).Value = 44232.55611753472
}

# 2) Append 14 new rows (604-617) following the existing 14-row cycle
#    (Nombre / URL / Disponibilidad / Fecha), refreshed by the
#    02-05-2021 13:51 automated availability check.
$names = @("Odoo","Blackbox","PowerBI","Dropbox","Odoo","GEE","UtilidadesOdoo","Filtros Dashboard","MapStore","GeoServer","Tomcat","Shiny","Github","EZ Exporter")
# Displayed cell text (what shows in column B) - for the MapStore row this
# includes the "#/" fragment even though the hyperlink target itself does not.
$displayUrls  = @(
    "https://www.dataintelligence-group.com/",
    "https://serviciodashboard.azurewebsites.net/",
    "https://powerbi.microsoft.com/es-es/",
    "https://www.dropbox.com/",
    "https://dataintelligence.store/",
    "https://app-data-i.users.earthengine.app/",
    "https://odooutil.azurewebsites.net/",
    "https://filtradordashboard.azurewebsites.net/",
    "https://ide.dataintelligence-group.com/mapstore/#/",
    "https://ide.dataintelligence-group.com/geoserver/web/?0",
    "https://ide.dataintelligence-group.com/",
    "https://rpubs.com/dataintelligence/",
    "https://github.com/Sud-Austral/",
    "https://ezexporter.highviewapps.com/exports/export-profile/"
)
# Hyperlink target (the relationship's Target) - base address only.
$linkUrls = @(
    "https://www.dataintelligence-group.com/",
    "https://serviciodashboard.azurewebsites.net/",
    "https://powerbi.microsoft.com/es-es/",
    "https://www.dropbox.com/",
    "https://dataintelligence.store/",
    "https://app-data-i.users.earthengine.app/",
    "https://odooutil.azurewebsites.net/",
    "https://filtradordashboard.azurewebsites.net/",
    "https://ide.dataintelligence-group.com/mapstore/",
    "https://ide.dataintelligence-group.com/geoserver/web/?0",
    "https://ide.dataintelligence-group.com/",
    "https://rpubs.com/dataintelligence/",
    "https://github.com/Sud-Austral/",
    "https://ezexporter.highviewapps.com/exports/export-profile/"
)
$timestamp = 44232.57720749198
$startRow = 604

for ($i = 0; $i -lt 14; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $names[$i]
    $ws.Cells.Item($r, 2).Value = $displayUrls[$i]
    $ws.Cells.Item($r, 3).Value = "Disponible"
    $ws.Cells.Item($r, 4).Value = $timestamp
    # Match the date/time formatting used by every prior row in column D.
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    if ($r -eq 612) {
        $ws.Hyperlinks.Add($ws.Cells.Item($r, 2), $linkUrls[$i], "/")
    } else {
        $ws.Hyperlinks.Add($ws.Cells.Item($r, 2), $linkUrls[$i])
    }

    # Match the hyperlink-cell look used by every prior row in column B
    # (Hyperlinks.Add stamps its own style, so re-apply the shared one).
    $ws.Cells.Item($r, 2).Style = "Hyperlink"
}
